$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the duplicate ID bug: B2 value changes from 3 to 101
$ws.Range("B2").Value = 101

# Add new rows of test data
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 3

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 101

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3

# Apply the same style as A2 to the new A column cells
$src = $ws.Range("A2")
$dst = $ws.Range("A3:A5")
$src.Copy()
$dst.PasteSpecial(-4122)
